$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.495.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.32%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.31%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.411.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  +4.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.001.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.25%  "

$ws.Range("E14").Value = "  -2.95%  "

$ws.Range("E15").Value = "  +8.66%  "

$ws.Range("E16").Value = "  +1.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.531.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.444.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.11%  "

$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  +2.57%  "

$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("E26").Value = "  +29.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.03%  "

$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("E30").Value = "  +8.84%  "

$ws.Range("E31").Value = "  +5.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.79%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  -1.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.976.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0761"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("E42").Value = "  -3.15%  "

$ws.Range("E43").Value = "  +1.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.80%  "

$ws.Range("E45").Value = "  +3.04%  "

$ws.Range("E46").Value = "  +2.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.77%  "

$ws.Range("E48").Value = "  +3.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +24.17%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "294.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.23%  "
